# Apply updated dSF ("F" column) values, as re-pulled/pushed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = -9
    "F5"  = -5
    "F6"  = -5
    "F7"  = -7
    "F13" = 1
    "F14" = 7
    "F16" = -4
    "F17" = 7
    "F20" = -1
    "F21" = 3
    "F23" = 8
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
